# Apply updated cryptocurrency price/volume figures (row-by-row) as captured in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ('Price') cells that look like plain decimals (e.g. '508.67') need a leading
# single-quote so Excel keeps them as text instead of converting them to numbers -
# matching the original inlineStr/text storage used throughout the sheet. Values that
# already contain thousands separators (e.g. '57.465.13') or other non-numeric characters
# are unambiguous and do not need the marker. Column E (percentages) is always text-safe.

$ws.Range('D2').Value = '57.465.13'
$ws.Range('E2').Value = '  +1.33%  '
$ws.Range('D3').Value = '3.014.29'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '''508.67'
$ws.Range('E5').Value = '  -0.65%  '
$ws.Range('D6').Value = '''139.71'
$ws.Range('E6').Value = '  +0.28%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('D9').Value = '''7.57'
$ws.Range('E9').Value = '  +0.21%  '
$ws.Range('E10').Value = '  +0.98%  '
$ws.Range('E11').Value = '  +2.39%  '
$ws.Range('D12').Value = '3.527.25'
$ws.Range('E13').Value = '  +0.27%  '
$ws.Range('D14').Value = '''26.38'
$ws.Range('E14').Value = '  +2.24%  '
$ws.Range('E15').Value = '  +3.17%  '
$ws.Range('D16').Value = '57.423.30'
$ws.Range('E16').Value = '  +1.17%  '
$ws.Range('D17').Value = '''6.20'
$ws.Range('E17').Value = '  +4.53%  '
$ws.Range('D18').Value = '3.010.32'
$ws.Range('E18').Value = '  +0.14%  '
$ws.Range('D19').Value = '''12.84'
$ws.Range('E19').Value = '  +2.11%  '
$ws.Range('D20').Value = '''7.96'
$ws.Range('E20').Value = '  +0.92%  '
$ws.Range('D21').Value = '''328.16'
$ws.Range('E21').Value = '  -1.04%  '
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('E23').Value = '  -1.73%  '
$ws.Range('D24').Value = '''0.500'
$ws.Range('D25').Value = '''64.48'
$ws.Range('E25').Value = '  +2.27%  '
$ws.Range('E26').Value = '  -3.16%  '
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').Value = '0.0₃0920'
$ws.Range('E28').Value = '  +1.13%  '
$ws.Range('E29').Value = '  +0.58%  '
$ws.Range('D30').Value = '''7.37'
$ws.Range('E30').Value = '  +3.58%  '
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('E32').Value = '  -5.44%  '
$ws.Range('D33').Value = '''20.58'
$ws.Range('E33').Value = '  -0.72%  '
$ws.Range('D34').Value = '''4.77'
$ws.Range('E34').Value = '  +3.57%  '
$ws.Range('D35').Value = '''153.89'
$ws.Range('E35').Value = '  -0.32%  '
$ws.Range('E36').Value = '  +3.25%  '
$ws.Range('D37').Value = '''1.28'
$ws.Range('E37').Value = '  -0.31%  '
$ws.Range('D38').Value = '''24.53'
$ws.Range('E38').Value = '  +2.87%  '
$ws.Range('D39').Value = '''0.0678'
$ws.Range('E39').Value = '  -0.33%  '
$ws.Range('D40').Value = '3.044.28'
$ws.Range('E40').Value = '  +0.15%  '
$ws.Range('D41').Value = '''37.83'
$ws.Range('E41').Value = '  +2.33%  '
$ws.Range('E42').Value = '  +4.59%  '
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').Value = '''0.650'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('E45').Value = '  -0.57%  '
$ws.Range('D46').Value = '2.222.55'
$ws.Range('E46').Value = '  -2.40%  '
$ws.Range('D47').Value = '''0.982'
$ws.Range('E47').Value = '  -2.42%  '
$ws.Range('D48').Value = '''6.06'
$ws.Range('E48').Value = '  +3.56%  '
$ws.Range('E49').Value = '  -0.85%  '
$ws.Range('D50').Value = '''19.56'
$ws.Range('E50').Value = '  -0.66%  '
$ws.Range('E51').Value = '  -5.58%  '
